# Fix: prevent hidden columns from being labeled upon detecting changes.
#
# The "Änderung" (change) marker in column L had been erroneously applied
# to every data row of the 55010 table (rows 52-90). In reality only the
# first row of each logical group should be rendered as a bold / shaded
# "group header" row (matching the pattern already used by rows 2-50),
# and the "ÄNDERUNG" label itself should not appear at all in this sheet.
#
# This script:
#   1. Re-applies the existing "group header" formatting (as already used
#      for rows 2-50) to the first row of every group in rows 52-90.
#   2. Clears the erroneous "ÄNDERUNG" text from column L for every
#      affected row and restores the plain (non-highlighted, centered)
#      formatting already used elsewhere in column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that start a new logical group (currently mis-styled like a normal
# detail row) - restore the bold/shaded "group header" look used by row 2.
$groupHeaderRows = @(52, 56, 60, 63, 66, 70, 76, 79, 85, 88)

foreach ($r in $groupHeaderRows) {
    $ws.Range("A2:V2").Copy()
    $ws.Range("A$r" + ":V$r").PasteSpecial(-4122)
    $ws.Range("L$r").Value = ""
}

# Every data row in this block (52-90) has a spurious "ÄNDERUNG" flag in
# column L. Clear it and reset column L's formatting to match the plain
# "no change" look used elsewhere (e.g. L3), for every row that wasn't
# already just reformatted above.
$detailRows = @(53, 54, 55, 57, 58, 59, 61, 62, 64, 65, 67, 68, 69, 71, 72, 73, 74, 75, 77, 78, 80, 81, 82, 83, 84, 86, 87, 89, 90)

foreach ($r in $detailRows) {
    $ws.Range("L3").Copy()
    $ws.Range("L$r").PasteSpecial(-4122)
    $ws.Range("L$r").Value = ""
}

$excel.CutCopyMode = 0
